$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove old rows 3 and 4 (their content is superseded); this shifts old rows 5,6 up to 3,4
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()

# ---- Row 2 ----
$ws.Range("A2").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector', RandomUnderSampler(random_state=42)),`n                ('model',`n                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',`n                                                                    criterion='entropy',`n                                                                    max_depth=1,`n                                                                    max_features='sqrt',`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=6,`n                                                                    random_state=42),`n                                   n_estimators=5, random_state=42))])"
$ws.Range("B2").Value = 0.6452147852147851
$ws.Range("C2").Value = "{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 1, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D2").Value = 0.830674648951034
$ws.Range("E2").Value = 0.5442775280275279
$ws.Range("F2").Value = 0.742857142857143
$ws.Range("G2").Value = 0.8403420944307578
$ws.Range("H2").Value = 0.5508154761904762
$ws.Range("I2").Value = 0.6842105263157895
$ws.Range("J2").Value = 0.8325531914893618
$ws.Range("K2").Value = 0.5583333333333335
$ws.Range("L2").Value = 0.8125
$ws.Range("M2").Value = "[1 0 1 1 1 1 0 1 0 1 0 1 0 1 1 0 0 1 1 1 1 0 1 1]"
$ws.Range("N2").Value = "[1 0 1 1 1 1 1 1 1 0 1 1 1 1 1 0 1 1 0 1 1 1 1 0]"
$ws.Range("O2").Value = 42
$ws.Rows.Item(2).AutoFit()

# ---- Row 3 ----
$ws.Range("A3").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector', RandomUnderSampler(random_state=42)),`n                ('model',`n                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',`n                                                                    criterion='entropy',`n                                                                    max_depth=2,`n                                                                    max_features='sqrt',`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=6,`n                                                                    random_state=42),`n                                   n_estimators=5, random_state=42))])"
$ws.Range("B3").Value = 0.6304295704295704
$ws.Range("C3").Value = "{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D3").Value = 0.8100822248624148
$ws.Range("E3").Value = 0.536361777111777
$ws.Range("F3").Value = 0.6060606060606061
$ws.Range("G3").Value = 0.8399613600472022
$ws.Range("H3").Value = 0.6584285714285714
$ws.Range("I3").Value = 0.5882352941176471
$ws.Range("J3").Value = 0.795
$ws.Range("K3").Value = 0.4858333333333333
$ws.Range("L3").Value = 0.625
$ws.Range("M3").Value = "[1 1 0 1 0 0 1 0 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 0]"
$ws.Range("N3").Value = "[0 1 1 1 1 1 1 1 0 1 0 1 1 0 1 1 0 1 1 1 1 0 0 1]"
$ws.Range("O3").Value = 69
$ws.Rows.Item(3).AutoFit()

# ---- Row 4 ----
$ws.Range("A4").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector', RandomUnderSampler(random_state=42)),`n                ('model',`n                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',`n                                                                    criterion='entropy',`n                                                                    max_depth=2,`n                                                                    max_features='sqrt',`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=6,`n                                                                    random_state=42),`n                                   n_estimators=5, random_state=42))])"
$ws.Range("B4").Value = 0.6724708624708624
$ws.Range("C4").Value = "{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 5, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 2, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D4").Value = 0.8237827983682544
$ws.Range("E4").Value = 0.6012970640470641
$ws.Range("F4").Value = 0.8
$ws.Range("G4").Value = 0.8354723602136482
$ws.Range("H4").Value = 0.6038670634920635
$ws.Range("I4").Value = 0.7619047619047619
$ws.Range("J4").Value = 0.8245555555555555
$ws.Range("K4").Value = 0.6289999999999999
$ws.Range("L4").Value = 0.8421052631578947
$ws.Range("M4").Value = "[0 1 0 0 1 1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 1 1 1]"
$ws.Range("N4").Value = "[1 1 1 1 1 1 0 1 1 1 1 1 1 1 1 0 0 1 1 1 1 1 1 1]"
$ws.Range("O4").Value = 23
$ws.Rows.Item(4).AutoFit()

# ---- Row 5 ----
$ws.Range("A5").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector', RandomUnderSampler(random_state=42)),`n                ('model',`n                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',`n                                                                    criterion='entropy',`n                                                                    max_depth=1,`n                                                                    max_features='sqrt',`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=6,`n                                                                    random_state=42),`n                                   random_state=42))])"
$ws.Range("B5").Value = 0.6223748473748474
$ws.Range("C5").Value = "{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 10, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 1, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D5").Value = 0.8211886316546199
$ws.Range("E5").Value = 0.5393069985569985
$ws.Range("F5").Value = 0.4827586206896552
$ws.Range("G5").Value = 0.8583036109560018
$ws.Range("H5").Value = 0.5832003968253968
$ws.Range("I5").Value = 0.4666666666666667
$ws.Range("J5").Value = 0.7964285714285713
$ws.Range("K5").Value = 0.5349999999999999
$ws.Range("L5").Value = 0.5
$ws.Range("M5").Value = "[0 1 1 0 0 1 0 0 0 0 1 1 1 0 0 1 1 0 1 1 1 1 1 1]"
$ws.Range("N5").Value = "[1 1 0 1 0 0 1 1 0 1 1 0 0 1 1 0 0 1 0 1 1 1 1 1]"
$ws.Range("O5").Value = 99
$ws.Rows.Item(5).AutoFit()

# ---- Row 6 ----
$ws.Range("A6").Value = "Pipeline(steps=[('scaler', MinMaxScaler()),`n                ('selector', RandomUnderSampler(random_state=42)),`n                ('model',`n                 BaggingClassifier(estimator=DecisionTreeClassifier(class_weight='balanced',`n                                                                    criterion='entropy',`n                                                                    max_depth=1,`n                                                                    max_features='sqrt',`n                                                                    min_samples_leaf=4,`n                                                                    min_samples_split=6,`n                                                                    random_state=42),`n                                   n_estimators=200, random_state=42))])"
$ws.Range("B6").Value = 0.6349494949494948
$ws.Range("C6").Value = "{'selector': RandomUnderSampler(random_state=42), 'scaler': MinMaxScaler(), 'model__n_estimators': 200, 'model__estimator__min_samples_split': 6, 'model__estimator__min_samples_leaf': 4, 'model__estimator__max_features': 'sqrt', 'model__estimator__max_depth': 1, 'model__estimator__criterion': 'entropy', 'model__estimator__class_weight': 'balanced'}"
$ws.Range("D6").Value = 0.8107816526140612
$ws.Range("E6").Value = 0.5848952713952714
$ws.Range("F6").Value = 0.5185185185185185
$ws.Range("G6").Value = 0.8924376331610654
$ws.Range("H6").Value = 0.6626488095238094
$ws.Range("I6").Value = 0.4375
$ws.Range("J6").Value = 0.7451923076923078
$ws.Range("K6").Value = 0.5425
$ws.Range("L6").Value = 0.6363636363636364
$ws.Range("M6").Value = "[1 0 1 1 0 0 0 0 1 0 1 1 0 1 1 0 1 0 0 0 0 0 1 1]"
$ws.Range("N6").Value = "[1 1 1 1 1 1 1 0 0 0 0 0 1 1 1 0 0 1 0 1 1 1 1 1]"
$ws.Range("O6").Value = 89
$ws.Rows.Item(6).AutoFit()
